$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new attendance week column (J) was added/filled in, largely mirroring
# the existing "24/9/2022" week recorded in column I, but with a handful
# of students marked absent (0) for this particular week.
$ws.Range("J4").Value = "24/9/2022"
$ws.Range("J5").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 0
$ws.Range("J8").Value = 3
$ws.Range("J9").Value = 3
$ws.Range("J10").Value = 3
$ws.Range("J11").Value = 3
$ws.Range("J12").Value = 3
$ws.Range("J13").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("J15").Value = 3
$ws.Range("J16").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("J19").Value = 3
$ws.Range("J20").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("J34").Value = 3
$ws.Range("J35").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("J41").Value = 0

# Reflect the author's final on-screen selection when the file was saved.
[void]$ws.Range("I44").Select()
